$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2026-02-25T18:07:33"
$ws.Range("V4").Value = 83.28
$ws.Range("W4").Value = 70.48
$ws.Range("X4").Value = 53.14
$ws.Range("Y4").Value = 57.06
$ws.Range("Z4").Value = 44.98
$ws.Range("V6").Value = -3.25
$ws.Range("W6").Value = -2.82
$ws.Range("X6").Value = -1.7
$ws.Range("Y6").Value = -1.26
$ws.Range("Z6").Value = -0.54
$ws.Range("V9").Value = 83.59999999999999
$ws.Range("W9").Value = 70.14
$ws.Range("X9").Value = 55
$ws.Range("Y9").Value = 58.72
$ws.Range("Z9").Value = 46.59
$ws.Range("V11").Value = -2.93
$ws.Range("W11").Value = -3.16
$ws.Range("X11").Value = 0.17
$ws.Range("Y11").Value = 0.41
$ws.Range("Z11").Value = 1.07
$ws.Range("V14").Value = 83.69
$ws.Range("W14").Value = 70.14
$ws.Range("X14").Value = 55
$ws.Range("Y14").Value = 58.78
$ws.Range("Z14").Value = 46.64
$ws.Range("V16").Value = -2.85
$ws.Range("W16").Value = -3.16
$ws.Range("X16").Value = 0.17
$ws.Range("Y16").Value = 0.47
$ws.Range("Z16").Value = 1.12
$ws.Range("V19").Value = 82.09999999999999
$ws.Range("W19").Value = 69.48
$ws.Range("X19").Value = 52.53
$ws.Range("Y19").Value = 56.45
$ws.Range("Z19").Value = 44.58
$ws.Range("V21").Value = -4.43
$ws.Range("W21").Value = -3.82
$ws.Range("X21").Value = -2.31
$ws.Range("Y21").Value = -1.86
$ws.Range("Z21").Value = -0.9399999999999999
$ws.Range("V24").Value = 82.09999999999999
$ws.Range("W24").Value = 69.48
$ws.Range("X24").Value = 52.53
$ws.Range("Y24").Value = 56.45
$ws.Range("Z24").Value = 44.58
$ws.Range("V26").Value = -4.43
$ws.Range("W26").Value = -3.82
$ws.Range("X26").Value = -2.31
$ws.Range("Y26").Value = -1.86
$ws.Range("Z26").Value = -0.9399999999999999
$ws.Range("V29").Value = 80.64
$ws.Range("W29").Value = 68.31
$ws.Range("X29").Value = 51.64
$ws.Range("Y29").Value = 55.48
$ws.Range("Z29").Value = 43.94
$ws.Range("V31").Value = -5.89
$ws.Range("W31").Value = -4.99
$ws.Range("X31").Value = -3.2
$ws.Range("Y31").Value = -2.83
$ws.Range("Z31").Value = -1.58
$ws.Range("V34").Value = 85
$ws.Range("W34").Value = 70.69
$ws.Range("X34").Value = 55.51
$ws.Range("Y34").Value = 59.2
$ws.Range("Z34").Value = 47.27
$ws.Range("V36").Value = -1.53
$ws.Range("W36").Value = -2.62
$ws.Range("X36").Value = 0.67
$ws.Range("Y36").Value = 0.89
$ws.Range("Z36").Value = 1.75
$ws.Range("V39").Value = 83.28
$ws.Range("W39").Value = 70.48
$ws.Range("X39").Value = 53.14
$ws.Range("Y39").Value = 57.06
$ws.Range("Z39").Value = 44.98
$ws.Range("V41").Value = -3.25
$ws.Range("W41").Value = -2.82
$ws.Range("X41").Value = -1.7
$ws.Range("Y41").Value = -1.26
$ws.Range("Z41").Value = -0.54
$ws.Range("V44").Value = 87.76000000000001
$ws.Range("W44").Value = 74.04000000000001
$ws.Range("X44").Value = 55.45
$ws.Range("Y44").Value = 59.08
$ws.Range("Z44").Value = 46.4
$ws.Range("V46").Value = 1.23
$ws.Range("W46").Value = 0.74
$ws.Range("X46").Value = 0.61
$ws.Range("Y46").Value = 0.77
$ws.Range("Z46").Value = 0.88
$ws.Range("V49").Value = 84.34
$ws.Range("W49").Value = 71.09999999999999
$ws.Range("X49").Value = 55.17
$ws.Range("Y49").Value = 59.5
$ws.Range("Z49").Value = 46.4
$ws.Range("V51").Value = -2.19
$ws.Range("W51").Value = -2.2
$ws.Range("X51").Value = 0.33
$ws.Range("Y51").Value = 1.19
$ws.Range("Z51").Value = 0.88
$ws.Range("V54").Value = 86.7
$ws.Range("W54").Value = 73.97
$ws.Range("X54").Value = 54.95
$ws.Range("Y54").Value = 58.31
$ws.Range("Z54").Value = 45.89
$ws.Range("V56").Value = 0.17
$ws.Range("W56").Value = 0.67
$ws.Range("Z56").Value = 0.37
$ws.Range("V59").Value = 90.23
$ws.Range("W59").Value = 76.2
$ws.Range("X59").Value = 57.01
$ws.Range("Y59").Value = 60.62
$ws.Range("Z59").Value = 47.46
$ws.Range("V61").Value = 3.7
$ws.Range("W61").Value = 2.9
$ws.Range("X61").Value = 2.17
$ws.Range("Y61").Value = 2.3
$ws.Range("Z61").Value = 1.95
$ws.Range("V64").Value = 91.95999999999999
$ws.Range("W64").Value = 77.56999999999999
$ws.Range("X64").Value = 58.03
$ws.Range("Y64").Value = 61.58
$ws.Range("Z64").Value = 48.22
$ws.Range("V66").Value = 5.43
$ws.Range("W66").Value = 4.27
$ws.Range("Z66").Value = 2.7
$ws.Range("V69").Value = 92.73999999999999
$ws.Range("W69").Value = 78.23
$ws.Range("X69").Value = 58.59
$ws.Range("Y69").Value = 62.23
$ws.Range("Z69").Value = 48.74
$ws.Range("V71").Value = 6.21
$ws.Range("W71").Value = 4.93
$ws.Range("X71").Value = 3.75
$ws.Range("Y71").Value = 3.92
$ws.Range("Z71").Value = 3.22
$ws.Range("V74").Value = 90.04000000000001
$ws.Range("W74").Value = 75.95999999999999
$ws.Range("X74").Value = 56.89
$ws.Range("Y74").Value = 60.43
$ws.Range("Z74").Value = 47.42
$ws.Range("V76").Value = 3.51
$ws.Range("W76").Value = 2.66
$ws.Range("X76").Value = 2.05
$ws.Range("Y76").Value = 2.11
$ws.Range("Z76").Value = 1.9
$ws.Range("V79").Value = 86.53
$ws.Range("W79").Value = 73.3
$ws.Range("X79").Value = 54.84
$ws.Range("Y79").Value = 58.31
$ws.Range("Z79").Value = 45.52
$ws.Range("V84").Value = 82.33
$ws.Range("W84").Value = 76.04000000000001
$ws.Range("X84").Value = 55.9
$ws.Range("Y84").Value = 59.14
$ws.Range("Z84").Value = 46.54
$ws.Range("V86").Value = -4.2
$ws.Range("W86").Value = 2.74
$ws.Range("X86").Value = 1.06
$ws.Range("Y86").Value = 0.83
$ws.Range("Z86").Value = 1.02
$ws.Range("V89").Value = 80.64
$ws.Range("W89").Value = 68.31
$ws.Range("X89").Value = 51.64
$ws.Range("Y89").Value = 55.48
$ws.Range("Z89").Value = 43.89
$ws.Range("V91").Value = -5.89
$ws.Range("W91").Value = -4.99
$ws.Range("X91").Value = -3.2
$ws.Range("Y91").Value = -2.83
$ws.Range("Z91").Value = -1.62

Write-Host "Applied all changes"